$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Insert a new row at position 19 (pushes old row 19.. down by one)
$ws.Rows.Item(19).Insert()

# Set the key value (A19) - short plain string
$ws.Range("A19").Value2 = "r109"

# Build the rich-text body for B19 using a here-string to keep quotes/newlines literal
$pegasusText = @'
<Bold>r109 Wild Pegasus</Bold>
<LineBreak/><LineBreak/>You encounter a wild Pegasus. Each character in your party is allowed one attempt to capture it.
<LineBreak/><LineBreak/>Roll one die. If the result is 5 plus, the character captures the Pegasus. You may add it as a winged mount to your party.
<LineBreak/><LineBreak/>
                                 <InlineUIContainer><Image Source='../../Images/Pegasus.gif' Height='300' Width='300'></Image></InlineUIContainer>
'@

$ws.Range("B19").Value2 = $pegasusText

# Colour the word "Pegasus" (the one right before ".gif") red, matching the other
# image-reference entries in this workbook (e.g. Falcon, Overcast, ...)
$redStart = $pegasusText.IndexOf("Pegasus.gif") + 1
$redLen = "Pegasus".Length
$redRange = $ws.Range("B19").Characters($redStart, $redLen)
$redRange.Font.Color = 255
$redRange.Font.Name = "Calibri"
$redRange.Font.Size = 11

# Restore default font formatting on the remainder of the text after "Pegasus"
$tailStart = $redStart + $redLen
$tailLen = $pegasusText.Length - ($tailStart - 1)
$tailRange = $ws.Range("B19").Characters($tailStart, $tailLen)
$tailRange.Font.Name = "Calibri"
$tailRange.Font.Size = 11
$tailRange.Font.Color = 0

# Row height for the new row
$ws.Rows.Item(19).RowHeight = 75

# Re-apply the sort state over the now-expanded range (A43:B217), matching the
# shift of the previously sorted block A42:B216 -> A43:B217
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A43:A217")) | Out-Null
$ws.Sort.SetRange($ws.Range("A43:B217"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()
